# Insert a new weekly price record at row 148 of the "Cilantro" price
# table, pushing the existing rows 148:208 down to 149:209 (the sheet's
# dimension grows from A1:R208 to A1:R209).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 148:208 down by one row. Excel's default row-insert copies
# the formatting (incl. the date number format on column D) from the row
# above, matching the s="2" style already present on the other D cells.
$ws.Rows("148").Insert()

# Populate the newly inserted row 148 with the new record. Columns
# A, B, C, E, F, G, H, I, R are identical for every row in this table.
$ws.Range("A148").Value = 4
$ws.Range("B148").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C148").Value = "Los Lagos"
$ws.Range("D148").Value = 44553
$ws.Range("E148").Value = 10
$ws.Range("F148").Value = 100112040
$ws.Range("G148").Value = "Cilantro"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 80
$ws.Range("K148").Value = 7000
$ws.Range("L148").Value = 7000
$ws.Range("M148").Value = 7000
$ws.Range("N148").Value = "$/docena de atados (2 kilos)"
$ws.Range("O148").Value = "Región de La Araucanía"
$ws.Range("P148").Value = 3500
$ws.Range("Q148").Value = 2
$ws.Range("R148").Value = "Hortaliza"
